$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H28").Value = 719.95654
$ws.Range("I28").Value = 702.4091
$ws.Range("K28").Value = 702.4091
$ws.Range("M28").Value = -217.4091
$ws.Range("H55").Value = 192.66667
$ws.Range("I55").Value = 163.66667
$ws.Range("J55").Value = 250.66667
$ws.Range("K55").Value = 163.66667
$ws.Range("L55").Value = 250.66667
$ws.Range("M55").Value = 50.33332999999999
$ws.Range("N55").Value = -678.6666700000001
$ws.Range("H62").Value = 1801
$ws.Range("I62").Value = 1052.5
$ws.Range("J62").Value = 2300
$ws.Range("K62").Value = 1052.5
$ws.Range("L62").Value = 2300
$ws.Range("M62").Value = -428.5
$ws.Range("N62").Value = -3548
$ws.Range("H65").Value = 1801
$ws.Range("I65").Value = 1052.5
$ws.Range("J65").Value = 2300
$ws.Range("K65").Value = 5262.5
$ws.Range("L65").Value = 11500
$ws.Range("M65").Value = -2142.5
$ws.Range("N65").Value = -17740
$ws.Range("H76").Value = 6750.375
$ws.Range("I76").Value = 6571.857
$ws.Range("K76").Value = 6571.857
$ws.Range("M76").Value = -6256.857
$ws.Range("H79").Value = 6750.375
$ws.Range("I79").Value = 6571.857
$ws.Range("K79").Value = 6571.857
$ws.Range("M79").Value = -5479.857
$ws.Range("H96").Value = 629
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 629
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 1887
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -4633
$ws.Range("H99").Value = 632
$ws.Range("I99").Value = 617.375
$ws.Range("J99").Value = 690.5
$ws.Range("K99").Value = 1852.125
$ws.Range("L99").Value = 2071.5
$ws.Range("M99").Value = -354.125
$ws.Range("N99").Value = -5067.5
$ws.Range("H100").Value = 2678.077
$ws.Range("I100").Value = 2531.5
$ws.Range("J100").Value = 3166.6667
$ws.Range("K100").Value = 2531.5
$ws.Range("L100").Value = 3166.6667
$ws.Range("M100").Value = -1990.5
$ws.Range("N100").Value = -4248.6667
$ws.Range("H127").Value = 1305.6428
$ws.Range("I127").Value = 340.66666
$ws.Range("J127").Value = 2029.375
$ws.Range("K127").Value = 1021.99998
$ws.Range("L127").Value = 6088.125
$ws.Range("M127").Value = 3938.00002
$ws.Range("N127").Value = -16008.125
$ws.Range("H132").Value = 2584.7585
$ws.Range("I132").Value = 2306.9524
$ws.Range("J132").Value = 3314
$ws.Range("K132").Value = 6920.8572
$ws.Range("L132").Value = 9942
$ws.Range("M132").Value = -4390.8572
$ws.Range("N132").Value = -15002
$ws.Range("H138").Value = 3695.8538
$ws.Range("I138").Value = 2146.8
$ws.Range("J138").Value = 4375.263
$ws.Range("K138").Value = 6440.400000000001
$ws.Range("L138").Value = 13125.789
$ws.Range("M138").Value = -1300.400000000001
$ws.Range("N138").Value = -23405.789

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 100
$ws.Range("K5").Value = 100
$ws.Range("M5").Value = 12
$ws.Range("H32").Value = 10671.307
$ws.Range("I32").Value = 10344.972
$ws.Range("J32").Value = 15240
$ws.Range("K32").Value = 10344.972
$ws.Range("L32").Value = 15240
$ws.Range("M32").Value = -10057.972
$ws.Range("N32").Value = -15814
$ws.Range("H36").Value = 14506.5
$ws.Range("I36").Value = 14506.5
$ws.Range("K36").Value = 14506.5
$ws.Range("M36").Value = -14160.5
$ws.Range("H45").Value = 1281.1428
$ws.Range("I45").Value = 1225.8462
$ws.Range("K45").Value = 1225.8462
$ws.Range("M45").Value = -848.8462
$ws.Range("H97").Value = 1340
$ws.Range("I97").Value = 1080
$ws.Range("K97").Value = 1080
$ws.Range("M97").Value = -584
$ws.Range("H102").Value = 2800.625
$ws.Range("I102").Value = 2557.8572
$ws.Range("K102").Value = 2557.8572
$ws.Range("M102").Value = -935.8571999999999
$ws.Range("H132").Value = 4323.3267
$ws.Range("I132").Value = 4500.811
$ws.Range("J132").Value = 3776.0833
$ws.Range("K132").Value = 13502.433
$ws.Range("L132").Value = 11328.2499
$ws.Range("M132").Value = -10972.433
$ws.Range("N132").Value = -16388.2499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("M4").Value = 15
$ws.Range("H105").Value = 8406472
$ws.Range("I105").Value = 12989865
$ws.Range("J105").Value = 3586.6667
$ws.Range("K105").Value = 12989865
$ws.Range("L105").Value = 3586.6667
$ws.Range("M105").Value = -12988118
$ws.Range("N105").Value = -7080.6667
$ws.Range("H134").Value = 3458.4736
$ws.Range("I134").Value = 3337.3635
$ws.Range("K134").Value = 10012.0905
$ws.Range("M134").Value = -7477.0905

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1974.7188
$ws.Range("J31").Value = 1523.9487
$ws.Range("L31").Value = 1523.9487
$ws.Range("N31").Value = -2113.9487
$ws.Range("H34").Value = 1974.7188
$ws.Range("J34").Value = 1523.9487
$ws.Range("L34").Value = 1523.9487
$ws.Range("N34").Value = -1927.9487
$ws.Range("H58").Value = 1685375.4
$ws.Range("I58").Value = 2471023.8
$ws.Range("K58").Value = 2471023.8
$ws.Range("M58").Value = -2470820.8
$ws.Range("H122").Value = 1500.375
$ws.Range("I122").Value = 1668.5
$ws.Range("J122").Value = 996
$ws.Range("K122").Value = 5005.5
$ws.Range("L122").Value = 2988
$ws.Range("M122").Value = -2555.5
$ws.Range("N122").Value = -7888
$ws.Range("H136").Value = 1685375.4
$ws.Range("I136").Value = 2471023.8
$ws.Range("K136").Value = 7413071.399999999
$ws.Range("M136").Value = -7410521.399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 216.66667
$ws.Range("I14").Value = 216.66667
$ws.Range("K14").Value = 650.00001
$ws.Range("M14").Value = -477.00001
$ws.Range("H68").Value = 1173.5977
$ws.Range("I68").Value = 841.60785
$ws.Range("J68").Value = 1643.9166
$ws.Range("K68").Value = 2524.82355
$ws.Range("L68").Value = 4931.7498
$ws.Range("M68").Value = -1713.82355
$ws.Range("N68").Value = -6553.7498
$ws.Range("H70").Value = 12117.4
$ws.Range("I70").Value = 17274.666
$ws.Range("J70").Value = 9907.143
$ws.Range("K70").Value = 51823.99800000001
$ws.Range("L70").Value = 29721.429
$ws.Range("M70").Value = -51508.99800000001
$ws.Range("N70").Value = -30351.429
$ws.Range("H71").Value = 1173.5977
$ws.Range("I71").Value = 841.60785
$ws.Range("J71").Value = 1643.9166
$ws.Range("K71").Value = 7574.47065
$ws.Range("L71").Value = 14795.2494
$ws.Range("M71").Value = -3518.47065
$ws.Range("N71").Value = -22907.2494
$ws.Range("H73").Value = 12117.4
$ws.Range("I73").Value = 17274.666
$ws.Range("J73").Value = 9907.143
$ws.Range("K73").Value = 51823.99800000001
$ws.Range("L73").Value = 29721.429
$ws.Range("M73").Value = -50731.99800000001
$ws.Range("N73").Value = -31905.429
$ws.Range("H75").Value = 8511
$ws.Range("J75").Value = 11151.875
$ws.Range("L75").Value = 33455.625
$ws.Range("N75").Value = -35451.625
$ws.Range("H78").Value = 8511
$ws.Range("J78").Value = 11151.875
$ws.Range("L78").Value = 100366.875
$ws.Range("N78").Value = -110350.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2793.6365
$ws.Range("I80").Value = 2550.5881
$ws.Range("J80").Value = 3620
$ws.Range("K80").Value = 2550.5881
$ws.Range("L80").Value = 3620
$ws.Range("M80").Value = -1552.5881
$ws.Range("N80").Value = -5616
$ws.Range("H83").Value = 2793.6365
$ws.Range("I83").Value = 2550.5881
$ws.Range("J83").Value = 3620
$ws.Range("K83").Value = 12752.9405
$ws.Range("L83").Value = 18100
$ws.Range("M83").Value = -7760.940500000001
$ws.Range("N83").Value = -28084

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 20471.908
$ws.Range("I61").Value = 35265.332
$ws.Range("J61").Value = 2719.8
$ws.Range("K61").Value = 35265.332
$ws.Range("L61").Value = 2719.8
$ws.Range("M61").Value = -35063.332
$ws.Range("N61").Value = -3123.8
$ws.Range("H96").Value = 49333.332
$ws.Range("J96").Value = 49333.332
$ws.Range("L96").Value = 49333.332
$ws.Range("N96").Value = -54825.332
$ws.Range("H113").Value = 20471.908
$ws.Range("I113").Value = 35265.332
$ws.Range("J113").Value = 2719.8
$ws.Range("K113").Value = 35265.332
$ws.Range("L113").Value = 2719.8
$ws.Range("M113").Value = -33095.332
$ws.Range("N113").Value = -7059.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 100000
$ws.Range("J15").Value = 100000
$ws.Range("L15").Value = 100000
$ws.Range("N15").Value = -100576
$ws.Range("H21").Value = 27000
$ws.Range("J21").Value = 18000
$ws.Range("L21").Value = 18000
$ws.Range("N21").Value = -18470
$ws.Range("H24").Value = 54005.8
$ws.Range("I24").Value = 50000
$ws.Range("J24").Value = 55007.25
$ws.Range("K24").Value = 50000
$ws.Range("L24").Value = 55007.25
$ws.Range("M24").Value = -49770
$ws.Range("N24").Value = -55467.25
$ws.Range("H35").Value = 27000
$ws.Range("J35").Value = 18000
$ws.Range("L35").Value = 18000
$ws.Range("N35").Value = -18580
$ws.Range("H81").Value = 82969.86
$ws.Range("I81").Value = 103779.63
$ws.Range("J81").Value = 6667.3335
$ws.Range("K81").Value = 207559.26
$ws.Range("L81").Value = 13334.667
$ws.Range("M81").Value = -206498.26
$ws.Range("N81").Value = -15456.667
$ws.Range("H84").Value = 82969.86
$ws.Range("I84").Value = 103779.63
$ws.Range("J84").Value = 6667.3335
$ws.Range("K84").Value = 1037796.3
$ws.Range("L84").Value = 66673.33499999999
$ws.Range("M84").Value = -1032492.3
$ws.Range("N84").Value = -77281.33499999999
$ws.Range("H101").Value = 125006550
$ws.Range("J101").Value = 125006550
$ws.Range("L101").Value = 125006550
$ws.Range("N101").Value = -125013040
$ws.Range("H113").Value = 411.47827
$ws.Range("I113").Value = 302.81818
$ws.Range("K113").Value = 908.45454
$ws.Range("M113").Value = 1261.54546
$ws.Range("H122").Value = 78126800
$ws.Range("I122").Value = 96155590
$ws.Range("K122").Value = 288466770
$ws.Range("M122").Value = -288464320

